# "problema giocatori e test altri classifier"
# Adds a 10th match-week row (week index 9) to the Calendar sheet and
# updates the sheet's active view/selection to reflect where the user
# was working after the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Add row 10: copy row 9's formatting for column A (thin left/right
#     border, no top/bottom - matches the "last row" look) then fill in
#     the week number + the 20 Opponent/IsHome pairs for week 9.
$ws.Range("A9").Copy($ws.Range("A10"))
$ws.Range("A10").Value = 9

$ws.Range("B10").Value = "VERONA"
$ws.Range("C10").Value = 1
$ws.Range("D10").Value = "JUVENTUS"
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = "CROTONE"
$ws.Range("G10").Value = 1
$ws.Range("H10").Value = "SPEZIA"
$ws.Range("I10").Value = 1
$ws.Range("J10").Value = "BOLOGNA"
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = "MILAN"
$ws.Range("M10").Value = 0
$ws.Range("N10").Value = "PARMA"
$ws.Range("O10").Value = 1
$ws.Range("P10").Value = "SASSUOLO"
$ws.Range("Q10").Value = 0
$ws.Range("R10").Value = "BENEVENTO"
$ws.Range("S10").Value = 0
$ws.Range("T10").Value = "UDINESE"
$ws.Range("U10").Value = 1
$ws.Range("V10").Value = "FIORENTINA"
$ws.Range("W10").Value = 1
$ws.Range("X10").Value = "ROMA"
$ws.Range("Y10").Value = 1
$ws.Range("Z10").Value = "GENOA"
$ws.Range("AA10").Value = 0
$ws.Range("AB10").Value = "NAPOLI"
$ws.Range("AC10").Value = 0
$ws.Range("AD10").Value = "TORINO"
$ws.Range("AE10").Value = 0
$ws.Range("AF10").Value = "INTER"
$ws.Range("AG10").Value = 1
$ws.Range("AH10").Value = "CAGLIARI"
$ws.Range("AI10").Value = 0
$ws.Range("AJ10").Value = "SAMPDORIA"
$ws.Range("AK10").Value = 1
$ws.Range("AL10").Value = "LAZIO"
$ws.Range("AM10").Value = 0
$ws.Range("AN10").Value = "ATALANTA"
$ws.Range("AO10").Value = 0

# --- Update the view: scroll so column C is the leftmost visible
#     column and the active selection sits on the newly added row.
$ws.Activate()
$ws.Range("C10").Select()
$excel.ActiveWindow.ScrollColumn = 3
$excel.ActiveWindow.ScrollRow = 1
